# Daily auto push: insert a new data point for 2026/02/23 (月) at row 865,
# shifting the existing rows 865:906 down to 866:907.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 865, pushing everything from 865 downward
# (rows 865-906) down by one row (to 866-907).
$ws.Rows(865).Insert()

# Populate the newly inserted row 865 with the new data point.
# The date column holds plain text like "2026/02/23" (not a real Excel date),
# so a leading apostrophe keeps it from being auto-converted into a date
# serial number; resetting the style afterwards keeps the cell unstyled,
# matching the rest of the data rows.
$ws.Range("A865").Value = "'2026/02/23"
$ws.Range("A865").Style = "Normal"
$ws.Range("B865").Value = "月"
$ws.Range("C865").Value = 7
$ws.Range("D865").Value = 201
